$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166, pushing the existing row 166 (and below) down to 167.
$ws.Rows.Item(166).Insert()

# Populate the new row 166 with the new weekly entry.
$ws.Cells.Item(166, 1).Value = 10
$ws.Cells.Item(166, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(166, 3).Value = "La Araucanía"
$ws.Cells.Item(166, 4).Value = 44753
$ws.Cells.Item(166, 4).NumberFormat = $ws.Cells.Item(167, 4).NumberFormat
$ws.Cells.Item(166, 5).Value = 9
$ws.Cells.Item(166, 6).Value = 100114007
$ws.Cells.Item(166, 7).Value = "Jengibre"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 80
$ws.Cells.Item(166, 11).Value = 20000
$ws.Cells.Item(166, 12).Value = 20000
$ws.Cells.Item(166, 13).Value = 20000
$ws.Cells.Item(166, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(166, 15).Value = "Perú"
$ws.Cells.Item(166, 16).Value = 1538
$ws.Cells.Item(166, 17).Value = 13
$ws.Cells.Item(166, 18).Value = "Hortaliza"
